# "Generate Report for Handback" - refresh the localization-status report
# after a handback event: the zh-cn target is back in sync (status +
# error cleared, handback timestamp bumped) and the de-de handback
# timestamp is bumped as well. Column widths are re-autofit to the new
# cell contents.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- zh-cn sheet: handback completed, in sync with en-US -------------------
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-30 22:53:57"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: new handback received, status/error mirror zh-cn's -------
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-30 22:54:11"
$wsDeDe.Range("P2").Value = ""

# --- Overview sheet: summary status cells mirror the same shared text ------
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- Column widths: re-autofit Status (wider text) / Error Detail (now
#     empty, shrinks back to fit the header) on both language sheets, and
#     the matching summary columns on Overview. ColumnWidth is in the
#     "characters" unit; the stored OOXML width is derived from it.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8333333333333

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8333333333333

$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667
